$wb = $excel.ActiveWorkbook
try { $excel.Application.EnableAutoComplete = $false } catch {}
try { $excel.Application.AutoCorrect.AutoFill = $false } catch {}
$ws1 = $wb.Worksheets.Item("TABLE_1")
$c = $ws1.Range("DQ13")
$c.Value = "11/01/2021"
Write-Output ("DQ13=" + $c.Value2)
